$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of column letters used in this edit to their 1-based column index
$colmap = @{
    "E" = 5; "G" = 7; "H" = 8; "I" = 9; "J" = 10; "K" = 11;
    "M" = 13; "N" = 14; "O" = 15; "P" = 16; "Q" = 17; "R" = 18; "S" = 19; "T" = 20
}

# New values recomputed by the NATMI pipeline re-run (per commit message),
# keyed by worksheet row number, then column letter.
$data = @{
    2 = @{ "E" = 3; "G" = 19.020839; "H" = 57.062517; "I" = 0.002925155182898989; "J" = 0.002925155182898989; "K" = 3; "M" = 64.44050866666667; "N" = 193.321526; "O" = 0.5284084541711456; "P" = 0.5284084541711456; "Q" = 1225.712540426771; "R" = 11031.41286384094; "S" = 0.00154567672840637; "T" = 0.00154567672840637 }
    3 = @{ "E" = 3; "G" = 19.020839; "H" = 57.062517; "I" = 0.002925155182898989; "J" = 0.002925155182898989; "K" = 3; "M" = 2.457277; "N" = 7.371831; "O" = 0.02014952966552172; "P" = 0.02014952966552173; "Q" = 46.73947019540299; "R" = 420.655231758627; "S" = 0.00005894050113407781; "T" = 0.00005894050113407782 }
    4 = @{ "E" = 3; "G" = 19.020839; "H" = 57.062517; "I" = 0.002925155182898989; "J" = 0.002925155182898989; "K" = 3; "M" = 19.791731; "N" = 59.375193; "O" = 0.162291052623097; "P" = 0.162291052623097; "Q" = 376.455328882309; "R" = 3388.097959940781; "S" = 0.0004747265137185849; "T" = 0.0004747265137185849 }
    5 = @{ "E" = 3; "G" = 19.020839; "H" = 57.062517; "I" = 0.002925155182898989; "J" = 0.002925155182898989; "K" = 3; "M" = 24.85018033333334; "N" = 74.55054100000001; "O" = 0.2037700453876647; "P" = 0.2037700453876647; "Q" = 472.6712792412997; "R" = 4254.041513171697; "S" = 0.0005960590043852898; "T" = 0.0005960590043852898 }
    6 = @{ "E" = 3; "G" = 19.020839; "H" = 57.062517; "I" = 0.002925155182898989; "J" = 0.002925155182898989; "K" = 3; "M" = 7.452280999999999; "N" = 22.356843; "O" = 0.06110827435625039; "P" = 0.06110827435625039; "Q" = 141.748637083759; "R" = 1275.737733753831; "S" = 0.0001787511854511992; "T" = 0.0001787511854511992 }
    7 = @{ "E" = 3; "G" = 19.020839; "H" = 57.062517; "I" = 0.002925155182898989; "J" = 0.002925155182898989; "K" = 3; "M" = 2.960099333333333; "N" = 8.880298; "O" = 0.02427264379632052; "P" = 0.02427264379632051; "Q" = 56.30357284334067; "R" = 506.732155590066; "S" = 0.00007100124980346796; "T" = 0.00007100124980346796 }
    8 = @{ "E" = 3; "G" = 3.243298333333333; "H" = 9.729894999999999; "I" = 0.0004987766801158274; "J" = 0.0004987766801158274; "K" = 3; "M" = 64.44050866666667; "N" = 193.321526; "O" = 0.5284084541711456; "P" = 0.5284084541711456; "Q" = 208.9997943577522; "R" = 1880.99814921977; "S" = 0.0002635578145166203; "T" = 0.0002635578145166203 }
    9 = @{ "E" = 3; "G" = 3.243298333333333; "H" = 9.729894999999999; "I" = 0.0004987766801158274; "J" = 0.0004987766801158274; "K" = 3; "M" = 2.457277; "N" = 7.371831; "O" = 0.02014952966552172; "P" = 0.02014952966552173; "Q" = 7.969682398638333; "R" = 71.727141587745; "S" = 0.0000100501155124643; "T" = 0.00001005011551246431 }
    10 = @{ "E" = 3; "G" = 3.243298333333333; "H" = 9.729894999999999; "I" = 0.0004987766801158274; "J" = 0.0004987766801158274; "K" = 3; "M" = 19.791731; "N" = 59.375193; "O" = 0.162291052623097; "P" = 0.162291052623097; "Q" = 64.19048816608166; "R" = 577.714393494735; "S" = 0.00008094699243985139; "T" = 0.00008094699243985139 }
    11 = @{ "E" = 3; "G" = 3.243298333333333; "H" = 9.729894999999999; "I" = 0.0004987766801158274; "J" = 0.0004987766801158274; "K" = 3; "M" = 24.85018033333334; "N" = 74.55054100000001; "O" = 0.2037700453876647; "P" = 0.2037700453876647; "Q" = 80.59654845813279; "R" = 725.368936123195; "S" = 0.0001016357467455109; "T" = 0.0001016357467455109 }
    12 = @{ "E" = 3; "G" = 3.243298333333333; "H" = 9.729894999999999; "I" = 0.0004987766801158274; "J" = 0.0004987766801158274; "K" = 3; "M" = 7.452280999999999; "N" = 22.356843; "O" = 0.06110827435625039; "P" = 0.06110827435625039; "Q" = 24.16997054683166; "R" = 217.529734921485; "S" = 0.00003047938221101772; "T" = 0.00003047938221101772 }
    13 = @{ "E" = 3; "G" = 3.243298333333333; "H" = 9.729894999999999; "I" = 0.0004987766801158274; "J" = 0.0004987766801158274; "K" = 3; "M" = 2.960099333333333; "N" = 8.880298; "O" = 0.02427264379632052; "P" = 0.02427264379632051; "Q" = 9.600485234301111; "R" = 86.40436710870999; "S" = 0.00001210662869036278; "T" = 0.00001210662869036278 }
    14 = @{ "E" = 3; "G" = 2517.581258333333; "H" = 7552.743774999999; "I" = 0.3871709269072258; "J" = 0.3871709269072258; "K" = 3; "M" = 64.44050866666667; "N" = 193.321526; "O" = 0.5284084541711456; "P" = 0.5284084541711456; "Q" = 162234.2168966667; "R" = 1460107.952070001; "S" = 0.2045843909870567; "T" = 0.2045843909870567 }
    15 = @{ "E" = 3; "G" = 2517.581258333333; "H" = 7552.743774999999; "I" = 0.3871709269072258; "J" = 0.3871709269072258; "K" = 3; "M" = 2.457277; "N" = 7.371831; "O" = 0.02014952966552172; "P" = 0.02014952966552173; "Q" = 6186.394521733558; "R" = 55677.55069560202; "S" = 0.007801312077344688; "T" = 0.00780131207734469 }
    16 = @{ "E" = 3; "G" = 2517.581258333333; "H" = 7552.743774999999; "I" = 0.3871709269072258; "J" = 0.3871709269072258; "K" = 3; "M" = 19.791731; "N" = 59.375193; "O" = 0.162291052623097; "P" = 0.162291052623097; "Q" = 49827.29103557485; "R" = 448445.6193201735; "S" = 0.06283437727283384; "T" = 0.06283437727283384 }
    17 = @{ "E" = 3; "G" = 2517.581258333333; "H" = 7552.743774999999; "I" = 0.3871709269072258; "J" = 0.3871709269072258; "K" = 3; "M" = 24.85018033333334; "N" = 74.55054100000001; "O" = 0.2037700453876647; "P" = 0.2037700453876647; "Q" = 62562.34827340359; "R" = 563061.1344606322; "S" = 0.07889383734866963; "T" = 0.07889383734866963 }
    18 = @{ "E" = 3; "G" = 2517.581258333333; "H" = 7552.743774999999; "I" = 0.3871709269072258; "J" = 0.3871709269072258; "K" = 3; "M" = 7.452280999999999; "N" = 22.356843; "O" = 0.06110827435625039; "P" = 0.06110827435625039; "Q" = 18761.72297743359; "R" = 168855.5067969023; "S" = 0.02365934722421052; "T" = 0.02365934722421052 }
    19 = @{ "E" = 3; "G" = 2517.581258333333; "H" = 7552.743774999999; "I" = 0.3871709269072258; "J" = 0.3871709269072258; "K" = 3; "M" = 2.960099333333333; "N" = 8.880298; "O" = 0.02427264379632052; "P" = 0.02427264379632051; "Q" = 7452.290604404994; "R" = 67070.61543964494; "S" = 0.009397661997110337; "T" = 0.009397661997110337 }
    20 = @{ "E" = 3; "G" = 3661.277099333334; "H" = 10983.831298; "I" = 0.5630563238111778; "J" = 0.5630563238111778; "K" = 3; "M" = 64.44050866666667; "N" = 193.321526; "O" = 0.5284084541711456; "P" = 0.5284084541711456; "Q" = 235934.5586506579; "R" = 2123411.027855921; "S" = 0.2975237216763524; "T" = 0.2975237216763524 }
    21 = @{ "E" = 3; "G" = 3661.277099333334; "H" = 10983.831298; "I" = 0.5630563238111778; "J" = 0.5630563238111778; "K" = 3; "M" = 2.457277; "N" = 7.371831; "O" = 0.02014952966552172; "P" = 0.02014952966552173; "Q" = 8996.772006818515; "R" = 80970.94806136665; "S" = 0.01134532009999293; "T" = 0.01134532009999293 }
    22 = @{ "E" = 3; "G" = 3661.277099333334; "H" = 10983.831298; "I" = 0.5630563238111778; "J" = 0.5630563238111778; "K" = 3; "M" = 19.791731; "N" = 59.375193; "O" = 0.162291052623097; "P" = 0.162291052623097; "Q" = 72463.01146646563; "R" = 652167.1031981906; "S" = 0.09137900347740742; "T" = 0.09137900347740742 }
    23 = @{ "E" = 3; "G" = 3661.277099333334; "H" = 10983.831298; "I" = 0.5630563238111778; "J" = 0.5630563238111778; "K" = 3; "M" = 24.85018033333334; "N" = 74.55054100000001; "O" = 0.2037700453876647; "P" = 0.2037700453876647; "Q" = 90983.39616873694; "R" = 818850.5655186324; "S" = 0.1147340126588153; "T" = 0.1147340126588153 }
    24 = @{ "E" = 3; "G" = 3661.277099333334; "H" = 10983.831298; "I" = 0.5630563238111778; "J" = 0.5630563238111778; "K" = 3; "M" = 7.452280999999999; "N" = 22.356843; "O" = 0.06110827435625039; "P" = 0.06110827435625039; "Q" = 27284.86576309691; "R" = 245563.7918678722; "S" = 0.03440740031347522; "T" = 0.03440740031347522 }
    25 = @{ "E" = 3; "G" = 3661.277099333334; "H" = 10983.831298; "I" = 0.5630563238111778; "J" = 0.5630563238111778; "K" = 3; "M" = 2.960099333333333; "N" = 8.880298; "O" = 0.02427264379632052; "P" = 0.02427264379632051; "Q" = 10837.7439008852; "R" = 97539.6951079668; "S" = 0.01366686558513442; "T" = 0.01366686558513442 }
    26 = @{ "E" = 3; "G" = 274.3202006666667; "H" = 822.960602; "I" = 0.04218684342756861; "J" = 0.04218684342756862; "K" = 3; "M" = 64.44050866666667; "N" = 193.321526; "O" = 0.5284084541711456; "P" = 0.5284084541711456; "Q" = 17677.33326850208; "R" = 159095.9994165187; "S" = 0.02229188472192168; "T" = 0.02229188472192169 }
    27 = @{ "E" = 3; "G" = 274.3202006666667; "H" = 822.960602; "I" = 0.04218684342756861; "J" = 0.04218684342756862; "K" = 3; "M" = 2.457277; "N" = 7.371831; "O" = 0.02014952966552172; "P" = 0.02014952966552173; "Q" = 674.0807197335846; "R" = 6066.726477602262; "S" = 0.000850045053138514; "T" = 0.0008500450531385142 }
    28 = @{ "E" = 3; "G" = 274.3202006666667; "H" = 822.960602; "I" = 0.04218684342756861; "J" = 0.04218684342756862; "K" = 3; "M" = 19.791731; "N" = 59.375193; "O" = 0.162291052623097; "P" = 0.162291052623097; "Q" = 5429.271619460688; "R" = 48863.44457514619; "S" = 0.006846547226705893; "T" = 0.006846547226705894 }
    29 = @{ "E" = 3; "G" = 274.3202006666667; "H" = 822.960602; "I" = 0.04218684342756861; "J" = 0.04218684342756862; "K" = 3; "M" = 24.85018033333334; "N" = 74.55054100000001; "O" = 0.2037700453876647; "P" = 0.2037700453876647; "Q" = 6816.906455642855; "R" = 61352.15810078569; "S" = 0.008596414999997963; "T" = 0.008596414999997965 }
    30 = @{ "E" = 3; "G" = 274.3202006666667; "H" = 822.960602; "I" = 0.04218684342756861; "J" = 0.04218684342756862; "K" = 3; "M" = 7.452280999999999; "N" = 22.356843; "O" = 0.06110827435625039; "P" = 0.06110827435625039; "Q" = 2044.311219344387; "R" = 18398.80097409949; "S" = 0.002577965202396042; "T" = 0.002577965202396042 }
    31 = @{ "E" = 3; "G" = 274.3202006666667; "H" = 822.960602; "I" = 0.04218684342756861; "J" = 0.04218684342756862; "K" = 3; "M" = 2.960099333333333; "N" = 8.880298; "O" = 0.02427264379632052; "P" = 0.02427264379632051; "Q" = 812.0150431132662; "R" = 7308.135388019396; "S" = 0.001023986223408518; "T" = 0.001023986223408518 }
    32 = @{ "E" = 3; "G" = 27.06326066666667; "H" = 81.18978200000001; "I" = 0.004161973991013034; "J" = 0.004161973991013034; "K" = 3; "M" = 64.44050866666667; "N" = 193.321526; "O" = 0.5284084541711456; "P" = 0.5284084541711456; "Q" = 1743.970283538593; "R" = 15695.73255184733; "S" = 0.00219922224289171; "T" = 0.00219922224289171 }
    33 = @{ "E" = 3; "G" = 27.06326066666667; "H" = 81.18978200000001; "I" = 0.004161973991013034; "J" = 0.004161973991013034; "K" = 3; "M" = 2.457277; "N" = 7.371831; "O" = 0.02014952966552172; "P" = 0.02014952966552173; "Q" = 66.50192798120467; "R" = 598.5173518308421; "S" = 0.00008386181839904698; "T" = 0.00008386181839904699 }
    34 = @{ "E" = 3; "G" = 27.06326066666667; "H" = 81.18978200000001; "I" = 0.004161973991013034; "J" = 0.004161973991013034; "K" = 3; "M" = 19.791731; "N" = 59.375193; "O" = 0.162291052623097; "P" = 0.162291052623097; "Q" = 535.6287750975474; "R" = 4820.658975877926; "S" = 0.0006754511399914575; "T" = 0.0006754511399914575 }
    35 = @{ "E" = 3; "G" = 27.06326066666667; "H" = 81.18978200000001; "I" = 0.004161973991013034; "J" = 0.004161973991013034; "K" = 3; "M" = 24.85018033333334; "N" = 74.55054100000001; "O" = 0.2037700453876647; "P" = 0.2037700453876647; "Q" = 672.5269079746737; "R" = 6052.742171772064; "S" = 0.0008480856290510061; "T" = 0.0008480856290510061 }
    36 = @{ "E" = 3; "G" = 27.06326066666667; "H" = 81.18978200000001; "I" = 0.004161973991013034; "J" = 0.004161973991013034; "K" = 3; "M" = 7.452280999999999; "N" = 22.356843; "O" = 0.06110827435625039; "P" = 0.06110827435625039; "Q" = 201.6830232642473; "R" = 1815.147209378226; "S" = 0.0002543310485064029; "T" = 0.0002543310485064029 }
    37 = @{ "E" = 3; "G" = 27.06326066666667; "H" = 81.18978200000001; "I" = 0.004161973991013034; "J" = 0.004161973991013034; "K" = 3; "M" = 2.960099333333333; "N" = 8.880298; "O" = 0.02427264379632052; "P" = 0.02427264379632051; "Q" = 80.10993985722622; "R" = 720.9894587150361; "S" = 0.0001010221121734099; "T" = 0.0001010221121734099 }
}

foreach ($row in $data.Keys) {
    $rowdata = $data[$row]
    foreach ($col in $rowdata.Keys) {
        $colnum = $colmap[$col]
        $ws.Cells.Item($row, $colnum).Value = $rowdata[$col]
    }
}

